# Update of 12th April 2022
# Refresh the MDSi test result sheet with the newer job numbers for
# the FedEx/UPS/CPU locator rows (column B), keeping their values as
# text (matching the existing "Job#" shared-string column) and
# clearing the previously-applied border/box formatting from those
# three cells, same as in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "32341649"
$ws.Range("B2").ClearFormats()

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "32341650"
$ws.Range("B3").ClearFormats()

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "32341651"
$ws.Range("B4").ClearFormats()
